$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# 1. Update the BrowserVersion value used by existing rows (2-5), column F,
#    from "99.0.4844.51" to "94.0.4606.61". Using a leading apostrophe keeps
#    the quotePrefix style (s="1") intact on these cells.
$ws.Range("F2").Value = "'94.0.4606.61"
$ws.Range("F3").Value = "'94.0.4606.61"
$ws.Range("F4").Value = "'94.0.4606.61"
$ws.Range("F5").Value = "'94.0.4606.61"

# 2. Append two new rows (6 and 7) that mirror rows 4 and 5 respectively,
#    copying them so that cell styles / empty-string cells are preserved,
#    then updating the Browser/BrowserVersion columns to firefox / 92.0.1.
$ws.Range("A4:M4").Copy($ws.Range("A6:M6"))
$ws.Range("A5:M5").Copy($ws.Range("A7:M7"))

# Row 7 empty cells need to stay text-typed (shared empty string) like row 5.
$ws.Range("G7").Value = "'"
$ws.Range("H7").Value = "'"
$ws.Range("I7").Value = "'"
$ws.Range("J7").Value = "'"
$ws.Range("K7").Value = "'"
$ws.Range("L7").Value = "'"
$ws.Range("M7").Value = "'"

# Update Browser / BrowserVersion for the two new rows (F before E so that
# the new shared strings land in the same order as in the target workbook).
$ws.Range("F6").Value = "'92.0.1"
$ws.Range("E6").Value = "firefox"
$ws.Range("F7").Value = "'92.0.1"
$ws.Range("E7").Value = "firefox"

# 3. Update selection to C5 on the DATA sheet
$ws.Activate()
$ws.Range("C5").Select()
